$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "choices" sheet: add a new row (row 9) duplicating row 8's OOP-Payments
#    "forms" entry (MIF_V_OOP / OOP Payments Births / Pagamentos de Nascimento)
# ---------------------------------------------------------------------------
$wsChoices = $wb.Worksheets.Item("choices")
$wsChoices.Cells.Item(9, 1).Value = "forms"
$wsChoices.Cells.Item(9, 2).Value = "MIF_V_OOP"
$wsChoices.Cells.Item(9, 3).Value = "OOP Payments Births"
$wsChoices.Cells.Item(9, 4).Value = "Pagamentos de Nascimento"

# ---------------------------------------------------------------------------
# 2. "survey" sheet
# ---------------------------------------------------------------------------
$wsSurvey = $wb.Worksheets.Item("survey")

# 2a. Constraint correction: the OOP_ANC row's hash-string formula pointed at
#     the wrong variable name (MIF_V_ANC instead of MIF_V_OOP_ANC).
$wsSurvey.Cells.Item(30, 2).Value = "''?' + odkSurvey.getHashString('MIF_V_OOP_ANC')"

# 2b. New "OOP Vaccines" block (mirrors the OOP_ANC block directly above it)
$wsSurvey.Cells.Item(32, 1).Value = "MIF_V_OOP_VAC"

$wsSurvey.Cells.Item(33, 2).Value = "''?' + odkSurvey.getHashString('MIF_V_OOP_VAC')"
$wsSurvey.Cells.Item(33, 5).Value = "external_link"
$wsSurvey.Cells.Item(33, 7).Value = "Open form"

$wsSurvey.Cells.Item(34, 3).Value = "exit section"

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping, restored in the same order the
#    author last touched the sheets so that "initial" ends up active again.
# ---------------------------------------------------------------------------
$wsChoices.Activate()
$wsChoices.Range("B7").Select()

$wsSurvey.Activate()
$wsSurvey.Range("C34").Select()

$wsInitial = $wb.Worksheets.Item("initial")
$wsInitial.Activate()
$wsInitial.Range("A9").Select()
